$d = $word.ActiveDocument

# Locate the two adjacent "Lab Manager" work-experience paragraphs by their
# distinguishing text rather than a hard-coded index, so the script is
# resilient to any earlier structural differences in the document.
$koehlPara = $null
$edwardsPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($koehlPara -eq $null -and $t -like "*Lab Manager, Koehl Lab*") {
        $koehlPara = $p
    }
    elseif ($edwardsPara -eq $null -and $t -like "*Lab Manager, Edwards Lab*") {
        $edwardsPara = $p
    }
}

if ($koehlPara -eq $null -or $edwardsPara -eq $null) {
    throw "Could not locate the Koehl Lab / Edwards Lab work-experience paragraphs"
}

# Together they span the whole "Lab Manager" block that needs to become three
# paragraphs: the new current Ph.D. student entry, followed by the existing
# Koehl Lab and Edwards Lab entries (reordered, periods dropped from bullets).
$rng = $d.Range($koehlPara.Range.Start, $edwardsPara.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Ph.D. student</w:t></w:r><w:r><w:t xml:space="preserve">, Ingalls Lab, University of Washington, Seattle</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">2019 - 2025</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Published papers on marine microbial metabolomics with a focus on automated and untargeted methodologies</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Developed existing and novel mass spectrometry tools for analysis and visualization</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Fieldwork in the North Pacific Subtropical Gyre and California Current</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Mentored graduate and undergraduate students in metabolomics</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Lab Manager</w:t></w:r><w:r><w:t xml:space="preserve">, Koehl Lab, University of California, Berkeley</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">2018 - 2019</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">- Maintained live cultures of choanoflagellates and protozoa</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Managed and organized lab members, materials, and safety protocols</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Assisted visiting researchers with statistical analyses</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Tutored undergraduates in ImageJ and R programming</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Lab Manager</w:t></w:r><w:r><w:t xml:space="preserve">, Edwards Lab, University of California, Berkeley</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">2018 - 2019</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">- Planned and set up laboratory experiments</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Performed HPLC-MS sample preparation and analysis using Thermo Fusion Lumos ID-X</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Developed and maintained a lipidomics pipeline for data handling and statistical analysis</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">- Educated undergraduates in R programming and SLURM cluster usage</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng.InsertXML($xml)
Write-Output "Applied work experience restructuring"
